# Regenerate the schedule report with the new source data:
#  - the "Floor 3" sheet/room no longer exists and is removed entirely
#  - "Floor 1" now lists a single room (the second ITC 314/Francis Rivas
#    block is gone) with a refreshed course code ("ITC 131")
#  - "Floor 2" now shows a different class (room code "ITC" / teacher
#    "John  Doe") and becomes the active sheet/tab

$wb = $excel.ActiveWorkbook

# Drop the "Floor 3" worksheet (and its room/class data) completely.
$wb.Worksheets.Item("Floor 3").Delete()

# --- Floor 1 --------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Floor 1")

# Refresh the single remaining class row...
$ws1.Range("A6").Value = "ITC 131"
$ws1.Range("B6").Value = "Francis  Rivas"
$ws1.Range("C6").Value = "07:00:00"
$ws1.Range("D6").Value = "08:30:00"

# ...and remove the second room block that used to sit at I6:L6.
$ws1.Range("I6:L6").Clear()

# Selection on this sheet moves back to the remaining block.
$ws1.Range("A6:D6").Select()

# --- Floor 2 ----------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Floor 2")

$ws2.Range("A6").Value = "ITC"
$ws2.Range("B6").Value = "John  Doe"
$ws2.Range("C6").Value = "07:00:00"
$ws2.Range("D6").Value = "08:30:00"

$ws2.Range("A6:D6").Select()

# "Floor 2" is now the active/selected tab in the workbook.
$ws2.Activate()
